# Monitoreo a las actividades del 22 al 29 de abril, 2016
# Fills in rows 21 and 22 (items 18 and 19) of the "No Conformidades" report
# with the follow-up activities for April 29, 2016 (serial 42489), and moves
# the sheet's selection to F23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy number-format / alignment styles from the last filled-in row (17)
#     onto the D, E and G columns of rows 21-22, so they match the styling
#     used by the rest of the populated rows (date format on D/E, wrapped
#     left-aligned text on G) instead of the placeholder style.
$ws.Range("D17").Copy()
$ws.Range("D21:D22").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("E17").Copy()
$ws.Range("E21:E22").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("G17").Copy()
$ws.Range("G21:G22").PasteSpecial(-4122) # xlPasteFormats

# --- Row 21 (ID 18) ---
$ws.Range("B21").Value = "Algunas de las actividades retrasadas no cuentan con seguimiento."
$ws.Range("C21").Value = "Ventas"
$ws.Range("D21").Value = 42489
$ws.Range("E21").Value = 42489
$ws.Range("F21").Value = "Cerrada"
$ws.Range("G21").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."

# --- Row 22 (ID 19) ---
$ws.Range("B22").Value = "La actividad del 28 de abril no tiene comentarios"
$ws.Range("C22").Value = "Compras"
$ws.Range("D22").Value = 42489
$ws.Range("E22").Value = 42489
$ws.Range("F22").Value = "Cerrada"
$ws.Range("G22").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."

# --- Row heights to match the other populated rows ---
$ws.Rows.Item(21).RowHeight = 75
$ws.Rows.Item(22).RowHeight = 75

# --- Update the active selection / scroll position ---
$win = $excel.ActiveWindow
$ws.Range("F23").Select() | Out-Null
$win.ScrollRow = 19
$win.ScrollColumn = 1
